# Natmi following Dr Hou advice
# Refresh the Fgf7-Fgfr1 LR-pair sheet: a new sending cluster "sCs" is added
# to the analysis, which changes the pooled statistics for the existing
# "FAPs" sending-cluster rows and appends four new "sCs" sending-cluster rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("FAPs","Fgf7","Fgfr1","ECs",3,1,11.14069866666667,33.422096,0.9684602815214559,0.9684602815214559,3,1,3.76917,11.30751,0.02686645020528053,0.02686645020528053,41.99118719344001,377.92068474096,0.02601908992928816,0.02601908992928816),
    @("FAPs","Fgf7","Fgfr1","FAPs",3,1,11.14069866666667,33.422096,0.9684602815214559,0.9684602815214559,3,1,82.48060333333333,247.44181,0.5879175050094569,0.5879175050094571,918.8915475815289,8270.023928233761,0.5693747524128506,0.5693747524128507),
    @("FAPs","Fgf7","Fgfr1","M2",3,1,11.14069866666667,33.422096,0.9684602815214559,0.9684602815214559,3,1,0.5716056666666667,1.714817,0.004074375838860061,0.004074375838860062,6.368086488492445,57.31277839643201,0.003945871171926633,0.003945871171926634),
    @("FAPs","Fgf7","Fgfr1","sCs",3,1,11.14069866666667,33.422096,0.9684602815214559,0.9684602815214559,3,1,53.47143866666666,160.414316,0.3811416689464024,0.3811416689464024,595.7091854584818,5361.382669126336,0.3691205680073904,0.3691205680073905),
    @("sCs","Fgf7","Fgfr1","ECs",2,0.6666666666666666,0.3628176666666666,1.088453,0.03153971847854405,0.03153971847854405,3,1,3.76917,11.30751,0.02686645020528053,0.02686645020528053,1.36752146467,12.30769318203,0.0008473602759923698,0.0008473602759923698),
    @("sCs","Fgf7","Fgfr1","FAPs",2,0.6666666666666666,0.3628176666666666,1.088453,0.03153971847854405,0.03153971847854405,3,1,82.48060333333333,247.44181,0.5879175050094569,0.5879175050094571,29.92542004665889,269.32878041993,0.01854275259660628,0.01854275259660628),
    @("sCs","Fgf7","Fgfr1","M2",2,0.6666666666666666,0.3628176666666666,1.088453,0.03153971847854405,0.03153971847854405,3,1,0.5716056666666667,1.714817,0.004074375838860061,0.004074375838860062,0.2073886342334444,1.866497708101,0.0001285046669334281,0.0001285046669334281),
    @("sCs","Fgf7","Fgfr1","sCs",2,0.6666666666666666,0.3628176666666666,1.088453,0.03153971847854405,0.03153971847854405,3,1,53.47143866666666,160.414316,0.3811416689464024,0.3811416689464024,19.40038261034978,174.603443493148,0.01202110093901197,0.01202110093901197)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $colNum = $j + 1
        $ws.Cells.Item($rowNum, $colNum).Value = $rowValues[$j]
    }
}
